# S-401 - Feature Association Inland Features: add reviewer questions/notes
# in a new column G next to specific feature rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G10").Value = "Light?"
$ws.Range("G17").Value = "hours?"
$ws.Range("G19").Value = "connection to lock basin etc? (C_AGGR)"
$ws.Range("G33").Value = "composition? Restricted area anchorage area, berth, "
$ws.Range("G41").Value = "hours?"
$ws.Range("G47").Value = "are there movable structures over navigable water?"
$ws.Range("G56").Value = "notice marks"
$ws.Range("G60").Value = "notice marks"

# Restore the cursor position the author ended up at after typing the
# notes (selection resting below the data, in column G).
$ws.Range("G86").Select()
